$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2").Value = "08/16/2021"
$ws.Range("L2").Value = "57466607"
